$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 232, shifting existing rows 232:290 down to 233:290
$ws.Rows.Item(232).Insert()

# Fill the new row 232 with its values. Columns A,B,C,E,F,G,H,I,N,O,Q,R are
# identical to the surrounding rows (unchanged constants for this data block).
$ws.Cells.Item(232, 1).Value = 10
$ws.Cells.Item(232, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(232, 3).Value = "La Araucanía"
$ws.Cells.Item(232, 4).Value = 44841
$ws.Cells.Item(232, 4).NumberFormat = $ws.Cells.Item(233, 4).NumberFormat
$ws.Cells.Item(232, 5).Value = 9
$ws.Cells.Item(232, 6).Value = 100112039
$ws.Cells.Item(232, 7).Value = "Ciboulette"
$ws.Cells.Item(232, 8).Value = "Sin especificar"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 40
$ws.Cells.Item(232, 11).Value = 6000
$ws.Cells.Item(232, 12).Value = 7000
$ws.Cells.Item(232, 13).Value = 6500
$ws.Cells.Item(232, 14).Value = "$/docena de atados"
$ws.Cells.Item(232, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(232, 16).Value = 2167
$ws.Cells.Item(232, 17).Value = 3
$ws.Cells.Item(232, 18).Value = "Hortaliza"
